$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.126.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.143.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'567.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'149.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.09%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "3.135.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +16.99%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'36.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "3.648.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "65.158.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'543.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.28%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "3.142.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.99%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'79.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'26.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'554.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0451"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'52.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  +11.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "Kaspa"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "Maker"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "3.072.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.57%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.23%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'25.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "0.0₃0530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'119.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "  +3.68%  "
$ws.Range("E51").Style = "Normal"
